$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 126, shifting existing rows 126:212 down to 127:213
$ws.Rows("126:126").Insert()

# Populate the newly inserted row 126 with the new data record
$ws.Range("A126").Value = 10
$ws.Range("B126").Value = "Vega Modelo de Temuco"
$ws.Range("C126").Value = "La Araucanía"
$ws.Range("D126").Value = 45176
$ws.Range("E126").Value = 9
$ws.Range("F126").Value = 100112035
$ws.Range("G126").Value = "Bruselas (repollito)"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 150
$ws.Range("K126").Value = 25000
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = 25667
$ws.Range("N126").Value = "$/malla 15 kilos"
$ws.Range("O126").Value = "Región Metropolitana"
$ws.Range("P126").Value = 1711
$ws.Range("Q126").Value = 15
$ws.Range("R126").Value = "Hortaliza"
